# Renamed few transcripts. Updated the DataSheet.
# Column D (Speaker) values "RBD" -> "T" and "Student" -> "S" for the rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where Speaker was "RBD" -> now "T"
$rowsT = @(2,3,4,5,6,7,9,10,12,14,16,19,22,23,24,25,26,27,31,35,36,37,40,41,43,44,48,51,52,55,57)

# Rows where Speaker was "Student" -> now "S"
$rowsS = @(21,42)

foreach ($r in $rowsT) {
    $ws.Cells.Item($r, 4).Value = "T"
}

foreach ($r in $rowsS) {
    $ws.Cells.Item($r, 4).Value = "S"
}
